$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").Value = "30.298.75"
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").Value = "1.926.63"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'0.7475"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.71%  "

$ws.Range("D6").Value = "'242.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.33%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "'0.3153"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.33%  "

$ws.Range("D9").Value = "'27.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.30%  "

$ws.Range("D10").Value = "'0.06966"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.19%  "

$ws.Range("D11").Value = "'0.08003"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "'0.7690"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.86%  "

$ws.Range("D13").Value = "1.929.75"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "'5.319"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").Value = "'93.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.87%  "

$ws.Range("D16").Value = "'14.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.48%  "

$ws.Range("D17").Value = "30.292.48"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "'250.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.27%  "

$ws.Range("D19").Value = "'0.000007881"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.66%  "

$ws.Range("D20").Value = "'5.734"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.28%  "

$ws.Range("D21").Value = "2.184.26"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").Value = "'6.623"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.24%  "

$ws.Range("D25").Value = "'9.412"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.08%  "

$ws.Range("D26").Value = "'165.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("D27").Value = "'18.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.08%  "

$ws.Range("D28").Value = "'0.1318"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.99%  "

$ws.Range("D29").Value = "'2.180"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.78%  "

$ws.Range("D30").Value = "'1.371"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.82%  "

$ws.Range("D31").Value = "'1.506"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.42%  "

$ws.Range("D32").Value = "'4.365"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("D33").Value = "'4.093"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.73%  "

$ws.Range("D34").Value = "'0.05086"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.33%  "

$ws.Range("D35").Value = "'1.275"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "

$ws.Range("D36").Value = "'0.7424"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("D37").Value = "'2.774"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("D38").Value = "'0.01945"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.02%  "

$ws.Range("D39").Value = "'2.794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'76.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.59%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.381"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").Value = "'0.4422"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.36%  "

$ws.Range("D43").Value = "'1.947"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.95%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "'0.8304"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("D46").Value = "'100.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.37%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.680"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.34%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.429"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("D49").Value = "'37.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "

$ws.Range("D50").Value = "'970.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.99%  "

$ws.Range("E51").Value = "  -1.03%  "

